$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Update the cached "datetimeFigureOut" field text (6/3/2020 -> 6/4/2020)
#    on the slide master and every slide layout's Date Placeholder shape.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            if ($sh.TextFrame.TextRange.Text -eq "6/3/2020") {
                $sh.TextFrame.TextRange.Text = "6/4/2020"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $ly = $p.SlideMaster.CustomLayouts.Item($i)
    Update-DatePlaceholder $ly.Shapes
}

# ---------------------------------------------------------------------------
# 2) On slide 2: reorder "Straight Connector 33" so it sits immediately
#    before "Straight Connector 29" (i.e. move it to the very back of the
#    z-order), and delete "Straight Connector 50" entirely.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 33") {
        $sh.ZOrder(1)  # msoSendToBack
        break
    }
}

$toDelete = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "Straight Connector 50") {
        $toDelete = $sh
        break
    }
}
if ($toDelete -ne $null) {
    $toDelete.Delete()
}

# ---------------------------------------------------------------------------
# 3) Add four new label rectangles ("ALERTS"/"LOGS") by duplicating the
#    existing "Rectangle 44" shape (same no-fill / style / text formatting)
#    and repositioning + retexting each copy.
# ---------------------------------------------------------------------------
$template = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $sh = $s2.Shapes.Item($i)
    if ($sh.Name -eq "Rectangle 44") {
        $template = $sh
        break
    }
}

$EMU_PER_PT = 12700

$newLabels = @(
    @{ Name = "Rectangle 1";  X = 2985267;  Y = 2272472; CX = 1335341; CY = 243220; Text = "ALERTS" },
    @{ Name = "Rectangle 23"; X = 10730538; Y = 2271827; CX = 1335341; CY = 243220; Text = "LOGS" },
    @{ Name = "Rectangle 25"; X = 8308640;  Y = 2267029; CX = 1335341; CY = 243220; Text = "LOGS" },
    @{ Name = "Rectangle 43"; X = 6102351;  Y = 2267029; CX = 1196539; CY = 243220; Text = "ALERTS" }
)

foreach ($spec in $newLabels) {
    $dup = $template.Duplicate().Item(1)
    $dup.Name = $spec.Name
    $dup.Left = $spec.X / $EMU_PER_PT
    $dup.Top = $spec.Y / $EMU_PER_PT
    $dup.Width = $spec.CX / $EMU_PER_PT
    $dup.Height = $spec.CY / $EMU_PER_PT
    $dup.TextFrame.TextRange.Text = $spec.Text
    $dup.TextFrame.TextRange.Font.Size = 11
}
